# Update "广州-漫展信息.xlsx" to match the newer scrape (output generated at 456a3b4).
#
# Changes:
#  1. "展览" (sheet 1): bump "想去人数" (F column) counts for a batch of existing rows.
#  2. "演出" (sheet 2): same kind of F-column bumps, PLUS a brand-new event row
#     ("2024.04.24 广州·今泉爱夏 巡演") inserted before the last row (old row 13
#     shifts down to row 14, and its F value increments from 13 to 14 too).
#  3. "全部类型" (sheet 4): mirrors both of the above since it is the combined,
#     date-sorted listing of every sheet's events. The new event is inserted in
#     date order between the "2024.04.20" and "2024.04.28" rows (= row 39,
#     pushing the old row 39 down to row 40).
#  4. "本地生活" (sheet 3) is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "展览" - F column bumps only.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1Updates = @{
    2  = 117
    3  = 185
    4  = 427
    5  = 202
    7  = 1191
    8  = 394
    9  = 200
    10 = 55
    13 = 413
    14 = 793
    15 = 184
    17 = 291
    19 = 1019
    20 = 476
    21 = 272
    23 = 388
}
foreach ($row in $ws1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $ws1Updates[$row]
}

# ---------------------------------------------------------------------------
# 2) Sheet "演出" - F column bumps + new row insert.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2Updates = @{
    5  = 42
    6  = 44
    11 = 151
}
foreach ($row in $ws2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $ws2Updates[$row]
}

# Insert the new event as row 13; the former row 13 (夏川里美) shifts to row 14.
$ws2.Rows.Item(13).Insert()

# Carry the bold/centered/bordered "A column" formatting into the new row.
$ws2.Range("A14").Copy()
$ws2.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("A13").Value = 12
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2024.04.24"
$ws2.Range("C13").Value = "广州·今泉爱夏  巡演"
$ws2.Range("D13").Value = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$ws2.Range("E13").NumberFormat = "@"
$ws2.Range("E13").Value = "2024.04.24 20:00-04.24 21:30"
$ws2.Range("F13").Value = 3
$ws2.Range("G13").Value = 288
$ws2.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=81890"
$ws2.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"

# The shifted former-row-13 becomes row 14: its running index and "want to go"
# count both advance by one (12->13, 13->14).
$ws2.Range("A14").Value = 13
$ws2.Range("F14").Value = 14

# ---------------------------------------------------------------------------
# 3) Sheet "全部类型" - F column bumps (mirrors sheet 1 + sheet 2) + new row
#    insert at the date-sorted position (row 39).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4Updates = @{
    4  = 117
    5  = 185
    6  = 427
    7  = 202
    9  = 1191
    10 = 394
    11 = 200
    13 = 55
    16 = 42
    18 = 44
    20 = 413
    21 = 793
    22 = 184
    24 = 291
    26 = 1019
    27 = 476
    30 = 272
    32 = 388
    34 = 151
}
foreach ($row in $ws4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $ws4Updates[$row]
}

# Insert the new event as row 39 (between 2024.04.20 and 2024.04.28 in date
# order); the former row 39 (夏川里美) shifts to row 40.
$ws4.Rows.Item(39).Insert()

$ws4.Range("A40").Copy()
$ws4.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws4.Range("A39").Value = 38
$ws4.Range("B39").NumberFormat = "@"
$ws4.Range("B39").Value = "2024.04.24"
$ws4.Range("C39").Value = "广州·今泉爱夏  巡演"
$ws4.Range("D39").Value = "革新路124号太古仓码头54汇5号仓 太空间Livehouse"
$ws4.Range("E39").NumberFormat = "@"
$ws4.Range("E39").Value = "2024.04.24 20:00-04.24 21:30"
$ws4.Range("F39").Value = 3
$ws4.Range("G39").Value = 288
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=81890"
$ws4.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"

$ws4.Range("A40").Value = 39
$ws4.Range("F40").Value = 14

Write-Host "edit.ps1 complete"
